$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 2).Value = 1.270414358554376
$ws.Cells.Item(2, 3).Value = 0.1870054823450715
$ws.Cells.Item(2, 4).Value = 0.5166861284057234
$ws.Cells.Item(2, 5).Value = 0.1768352350388191
$ws.Cells.Item(2, 7).Value = 1.006974040763296
$ws.Cells.Item(2, 8).Value = 1.030110265252901
$ws.Cells.Item(2, 10).Value = 0.08044269317985986
$ws.Cells.Item(2, 12).Value = 0.3672628403549325
$ws.Cells.Item(2, 13).Value = 0.3377522682235536
$ws.Cells.Item(2, 14).Value = 1.682596187459559
$ws.Cells.Item(2, 15).Value = 4.128106805779993
# Row 3
$ws.Cells.Item(3, 2).Value = 1.194097904004337
$ws.Cells.Item(3, 3).Value = 0.1793620625467014
$ws.Cells.Item(3, 4).Value = 0.5168703469153542
$ws.Cells.Item(3, 5).Value = 0.1779761309630246
$ws.Cells.Item(3, 7).Value = 1.008550174120145
$ws.Cells.Item(3, 8).Value = 1.035315464045851
$ws.Cells.Item(3, 10).Value = 0.08013365584731602
$ws.Cells.Item(3, 12).Value = 0.36401133398887
$ws.Cells.Item(3, 13).Value = 0.3249708236499345
$ws.Cells.Item(3, 14).Value = 1.69762267941412
$ws.Cells.Item(3, 15).Value = 4.142023550315116
# Row 4
$ws.Cells.Item(4, 2).Value = 1.147566916431515
$ws.Cells.Item(4, 3).Value = 0.1746360996950784
$ws.Cells.Item(4, 4).Value = 0.5171931454736836
$ws.Cells.Item(4, 5).Value = 0.1787277619588679
$ws.Cells.Item(4, 7).Value = 1.010105835136784
$ws.Cells.Item(4, 8).Value = 1.038941806148443
$ws.Cells.Item(4, 10).Value = 0.07994477285580359
$ws.Cells.Item(4, 12).Value = 0.3621456192304322
$ws.Cells.Item(4, 13).Value = 0.317233421140628
$ws.Cells.Item(4, 14).Value = 1.707420632656948
$ws.Cells.Item(4, 15).Value = 4.152704089008267
# Row 5
$ws.Cells.Item(5, 2).Value = 1.128688933873775
$ws.Cells.Item(5, 3).Value = 0.172702091702817
$ws.Cells.Item(5, 4).Value = 0.5173775142707484
$ws.Cells.Item(5, 5).Value = 0.179046930732591
$ws.Cells.Item(5, 7).Value = 1.010887597981167
$ws.Cells.Item(5, 8).Value = 1.040527874795089
$ws.Cells.Item(5, 10).Value = 0.07986802785568514
$ws.Cells.Item(5, 12).Value = 0.3614183234236634
$ws.Cells.Item(5, 13).Value = 0.3141084038417929
$ws.Cells.Item(5, 14).Value = 1.711557222344069
$ws.Cells.Item(5, 15).Value = 4.157593666401539
# Row 6
$ws.Cells.Item(6, 2).Value = 1.12555935983562
$ws.Cells.Item(6, 3).Value = 0.1723804632471371
$ws.Cells.Item(6, 4).Value = 0.5174113220027152
$ws.Cells.Item(6, 5).Value = 0.1791007065361736
$ws.Cells.Item(6, 7).Value = 1.011026336514206
$ws.Cells.Item(6, 8).Value = 1.040797784434119
$ws.Cells.Item(6, 10).Value = 0.07985529832264859
$ws.Cells.Item(6, 12).Value = 0.3612995532438887
$ws.Cells.Item(6, 13).Value = 0.3135911981850938
$ws.Cells.Item(6, 14).Value = 1.712252790872917
$ws.Cells.Item(6, 15).Value = 4.158438025487271
# Row 7
$ws.Cells.Item(7, 2).Value = 1.147311980323167
$ws.Cells.Item(7, 3).Value = 0.1746100497776837
$ws.Cells.Item(7, 4).Value = 0.5171954179157439
$ws.Cells.Item(7, 5).Value = 0.1787320142321409
$ws.Cells.Item(7, 7).Value = 1.010115779775361
$ws.Cells.Item(7, 8).Value = 1.038962757777995
$ws.Cells.Item(7, 10).Value = 0.07994373691828294
$ws.Cells.Item(7, 12).Value = 0.3621356768883714
$ws.Cells.Item(7, 13).Value = 0.3171911621960248
$ws.Cells.Item(7, 14).Value = 1.70747583759298
$ws.Cells.Item(7, 15).Value = 4.152767856369366
# Row 8
$ws.Cells.Item(8, 2).Value = 1.244033254098724
$ws.Cells.Item(8, 3).Value = 0.1843769272719413
$ws.Cells.Item(8, 4).Value = 0.5167061840342058
$ws.Cells.Item(8, 5).Value = 0.1772180184725558
$ws.Cells.Item(8, 7).Value = 1.007395458919149
$ws.Cells.Item(8, 8).Value = 1.031815759195609
$ws.Cells.Item(8, 10).Value = 0.08033596199221371
$ws.Cells.Item(8, 12).Value = 0.3661146635643604
$ws.Cells.Item(8, 13).Value = 0.3333224534871846
$ws.Cells.Item(8, 14).Value = 1.687658743507733
$ws.Cells.Item(8, 15).Value = 4.132462142223005
# Row 9
$ws.Cells.Item(9, 2).Value = 1.436251943384605
$ws.Cells.Item(9, 3).Value = 0.2032642256542658
$ws.Cells.Item(9, 4).Value = 0.517406735393692
$ws.Cells.Item(9, 5).Value = 0.1746537908001731
$ws.Cells.Item(9, 7).Value = 1.006727666353598
$ws.Cells.Item(9, 8).Value = 1.021211218011885
$ws.Cells.Item(9, 10).Value = 0.08111166990549634
$ws.Cells.Item(9, 12).Value = 0.37494981231427
$ws.Cells.Item(9, 13).Value = 0.3658225886349911
$ws.Cells.Item(9, 14).Value = 1.653328232141455
$ws.Cells.Item(9, 15).Value = 4.109584364832813
# Row 10
$ws.Cells.Item(10, 2).Value = 1.578974954537443
$ws.Cells.Item(10, 3).Value = 0.2169738051612455
$ws.Cells.Item(10, 4).Value = 0.5189289051618005
$ws.Cells.Item(10, 5).Value = 0.1730153793699625
$ws.Cells.Item(10, 7).Value = 1.009086498808188
$ws.Cells.Item(10, 8).Value = 1.015494862560928
$ws.Cells.Item(10, 10).Value = 0.08168520023331638
$ws.Cells.Item(10, 12).Value = 0.3820648546817154
$ws.Cells.Item(10, 13).Value = 0.3902184490861629
$ws.Cells.Item(10, 14).Value = 1.63086074058026
$ws.Cells.Item(10, 15).Value = 4.103104847887465
# Row 11
$ws.Cells.Item(11, 2).Value = 1.644218564947266
$ws.Cells.Item(11, 3).Value = 0.22317340958665
$ws.Cells.Item(11, 4).Value = 0.5198392444761311
$ws.Cells.Item(11, 5).Value = 0.1723230846201886
$ws.Cells.Item(11, 7).Value = 1.010779408250556
$ws.Cells.Item(11, 8).Value = 1.013343982085289
$ws.Cells.Item(11, 10).Value = 0.08194682297250822
$ws.Cells.Item(11, 12).Value = 0.3854360384978435
$ws.Cells.Item(11, 13).Value = 0.4014272071624205
$ws.Cells.Item(11, 14).Value = 1.621236540773772
$ws.Cells.Item(11, 15).Value = 4.102400494066842
# Row 12
$ws.Cells.Item(12, 2).Value = 1.668969168930573
$ws.Cells.Item(12, 3).Value = 0.2255156128519218
$ws.Cells.Item(12, 4).Value = 0.5202152123542163
$ws.Cells.Item(12, 5).Value = 0.1720685365852663
$ws.Cells.Item(12, 7).Value = 1.011509660609789
$ws.Cells.Item(12, 8).Value = 1.012594060956602
$ws.Cells.Item(12, 10).Value = 0.08204598870258906
$ws.Cells.Item(12, 12).Value = 0.3867318483622597
$ws.Cells.Item(12, 13).Value = 0.40568739882697
$ws.Cells.Item(12, 14).Value = 1.617677780584984
$ws.Cells.Item(12, 15).Value = 4.102456306086395
# Row 13
$ws.Cells.Item(13, 2).Value = 1.663636738560967
$ws.Cells.Item(13, 3).Value = 0.2250114221726562
$ws.Cells.Item(13, 4).Value = 0.520132852892317
$ws.Cells.Item(13, 5).Value = 0.1721230199057917
$ws.Cells.Item(13, 7).Value = 1.011348420437471
$ws.Cells.Item(13, 8).Value = 1.012752699196199
$ws.Cells.Item(13, 10).Value = 0.08202462750850614
$ws.Cells.Item(13, 12).Value = 0.3864519199251646
$ws.Cells.Item(13, 13).Value = 0.4047691974476919
$ws.Cells.Item(13, 14).Value = 1.618440411683856
$ws.Cells.Item(13, 15).Value = 4.102429941684647
# Row 14
$ws.Cells.Item(14, 2).Value = 1.646253933439652
$ws.Cells.Item(14, 3).Value = 0.223366214126969
$ws.Cells.Item(14, 4).Value = 0.5198695499622659
$ws.Cells.Item(14, 5).Value = 0.1723019903978802
$ws.Cells.Item(14, 7).Value = 1.010837698890313
$ws.Cells.Item(14, 8).Value = 1.013280991927843
$ws.Cells.Item(14, 10).Value = 0.08195497954597997
$ws.Cells.Item(14, 12).Value = 0.3855422611728017
$ws.Cells.Item(14, 13).Value = 0.4017773830674258
$ws.Cells.Item(14, 14).Value = 1.620942041795701
$ws.Cells.Item(14, 15).Value = 4.102398621867934
# Row 15
$ws.Cells.Item(15, 2).Value = 1.635612183720241
$ws.Cells.Item(15, 3).Value = 0.2223577625213125
$ws.Cells.Item(15, 4).Value = 0.5197123352597259
$ws.Cells.Item(15, 5).Value = 0.1724126054234514
$ws.Cells.Item(15, 7).Value = 1.010536483051624
$ws.Cells.Item(15, 8).Value = 1.013612993425596
$ws.Cells.Item(15, 10).Value = 0.0819123302528908
$ws.Cells.Item(15, 12).Value = 0.3849875675970509
$ws.Cells.Item(15, 13).Value = 0.3999468439134048
$ws.Cells.Item(15, 14).Value = 1.622485523955966
$ws.Cells.Item(15, 15).Value = 4.102421440160896
# Row 16
$ws.Cells.Item(16, 2).Value = 1.574717406534546
$ws.Cells.Item(16, 3).Value = 0.216567890995492
$ws.Cells.Item(16, 4).Value = 0.5188737891929804
$ws.Cells.Item(16, 5).Value = 0.1730616881054647
$ws.Cells.Item(16, 7).Value = 1.008988340618544
$ws.Cells.Item(16, 8).Value = 1.015644467158779
$ws.Cells.Item(16, 10).Value = 0.08166811632631266
$ws.Cells.Item(16, 12).Value = 0.3818472357853295
$ws.Cells.Item(16, 13).Value = 0.3894881394298011
$ws.Cells.Item(16, 14).Value = 1.631501701523455
$ws.Cells.Item(16, 15).Value = 4.103196015182476
# Row 17
$ws.Cells.Item(17, 2).Value = 1.537440868951307
$ws.Cells.Item(17, 3).Value = 0.2130064292225029
$ws.Cells.Item(17, 4).Value = 0.5184151120268581
$ws.Cells.Item(17, 5).Value = 0.1734734491664778
$ws.Cells.Item(17, 7).Value = 1.008197402323376
$ws.Cells.Item(17, 8).Value = 1.017005793275899
$ws.Cells.Item(17, 10).Value = 0.08151847737332574
$ws.Cells.Item(17, 12).Value = 0.3799551057996808
$ws.Cells.Item(17, 13).Value = 0.3831002846279574
$ws.Cells.Item(17, 14).Value = 1.6371855487129
$ws.Cells.Item(17, 15).Value = 4.104245706910007
# Row 18
$ws.Cells.Item(18, 2).Value = 1.51603042967173
$ws.Cells.Item(18, 3).Value = 0.2109545007752445
$ws.Cells.Item(18, 4).Value = 0.5181718079629292
$ws.Cells.Item(18, 5).Value = 0.1737152753379458
$ws.Cells.Item(18, 7).Value = 1.007800817275069
$ws.Cells.Item(18, 8).Value = 1.017831106964167
$ws.Cells.Item(18, 10).Value = 0.08143247738596671
$ws.Cells.Item(18, 12).Value = 0.3788794715015911
$ws.Cells.Item(18, 13).Value = 0.3794366200281729
$ws.Cells.Item(18, 14).Value = 1.640510883174571
$ws.Cells.Item(18, 15).Value = 4.105060605588193
# Row 19
$ws.Cells.Item(19, 2).Value = 1.508786427564019
$ws.Cells.Item(19, 3).Value = 0.2102591613447657
$ws.Cells.Item(19, 4).Value = 0.5180929556198919
$ws.Cells.Item(19, 5).Value = 0.1737980114186382
$ws.Cells.Item(19, 7).Value = 1.007676559155357
$ws.Cells.Item(19, 8).Value = 1.018117813754301
$ws.Cells.Item(19, 10).Value = 0.08140337128009989
$ws.Cells.Item(19, 12).Value = 0.3785174599801735
$ws.Cells.Item(19, 13).Value = 0.3781979712133676
$ws.Cells.Item(19, 14).Value = 1.641646427047483
$ws.Cells.Item(19, 15).Value = 4.105372781505963
# Row 20
$ws.Cells.Item(20, 2).Value = 1.541405921878777
$ws.Cells.Item(20, 3).Value = 0.21338591302856
$ws.Cells.Item(20, 4).Value = 0.5184618164390855
$ws.Cells.Item(20, 5).Value = 0.1734290999529655
$ws.Cells.Item(20, 7).Value = 1.008275560549762
$ws.Cells.Item(20, 8).Value = 1.016856498855788
$ws.Cells.Item(20, 10).Value = 0.0815343996823259
$ws.Cells.Item(20, 12).Value = 0.3801552160148276
$ws.Cells.Item(20, 13).Value = 0.3837792020622857
$ws.Cells.Item(20, 14).Value = 1.636574683942115
$ws.Cells.Item(20, 15).Value = 4.104112113097074
# Row 21
$ws.Cells.Item(21, 2).Value = 1.651358495215504
$ws.Cells.Item(21, 3).Value = 0.2238496006901869
$ws.Cells.Item(21, 4).Value = 0.5199460412386543
$ws.Cells.Item(21, 5).Value = 0.172249216049785
$ws.Cells.Item(21, 7).Value = 1.010985289369145
$ws.Cells.Item(21, 8).Value = 1.013124067670091
$ws.Cells.Item(21, 10).Value = 0.08197543433567134
$ws.Cells.Item(21, 12).Value = 0.3858089295041225
$ws.Cells.Item(21, 13).Value = 0.4026557281801715
$ws.Cells.Item(21, 14).Value = 1.620204926418516
$ws.Cells.Item(21, 15).Value = 4.102399068021043
# Row 22
$ws.Cells.Item(22, 2).Value = 1.723476175814312
$ws.Cells.Item(22, 3).Value = 0.2306563806857866
$ws.Cells.Item(22, 4).Value = 0.5210981287958418
$ws.Cells.Item(22, 5).Value = 0.1715224371101822
$ws.Cells.Item(22, 7).Value = 1.013276110163787
$ws.Cells.Item(22, 8).Value = 1.011061043215577
$ws.Cells.Item(22, 10).Value = 0.0822642249301353
$ws.Cells.Item(22, 12).Value = 0.3896159030430226
$ws.Cells.Item(22, 13).Value = 0.4150838855066681
$ws.Cells.Item(22, 14).Value = 1.610005941732403
$ws.Cells.Item(22, 15).Value = 4.103159461126154
# Row 23
$ws.Cells.Item(23, 2).Value = 1.68496256013151
$ws.Cells.Item(23, 3).Value = 0.2270264336401056
$ws.Cells.Item(23, 4).Value = 0.5204666094930417
$ws.Cells.Item(23, 5).Value = 0.171906280619206
$ws.Cells.Item(23, 7).Value = 1.012005873244533
$ws.Cells.Item(23, 8).Value = 1.012127706369043
$ws.Cells.Item(23, 10).Value = 0.08211004462289395
$ws.Cells.Item(23, 12).Value = 0.3875738486876088
$ws.Cells.Item(23, 13).Value = 0.4084424856746409
$ws.Cells.Item(23, 14).Value = 1.615403631532324
$ws.Cells.Item(23, 15).Value = 4.102581619746701
# Row 24
$ws.Cells.Item(24, 2).Value = 1.539613257982296
$ws.Cells.Item(24, 3).Value = 0.2132143620846421
$ws.Cells.Item(24, 4).Value = 0.5184406378389212
$ws.Cells.Item(24, 5).Value = 0.1734491343475391
$ws.Cells.Item(24, 7).Value = 1.008240044129096
$ws.Cells.Item(24, 8).Value = 1.016923861944008
$ws.Cells.Item(24, 10).Value = 0.08152720111449341
$ws.Cells.Item(24, 12).Value = 0.3800647082489803
$ws.Cells.Item(24, 13).Value = 0.3834722360636178
$ws.Cells.Item(24, 14).Value = 1.636850676432104
$ws.Cells.Item(24, 15).Value = 4.104171852302244
# Row 25
$ws.Cells.Item(25, 2).Value = 1.383984417438228
$ws.Cells.Item(25, 3).Value = 0.198183660178529
$ws.Cells.Item(25, 4).Value = 0.5170399723173915
$ws.Cells.Item(25, 5).Value = 0.1753042791871344
$ws.Cells.Item(25, 7).Value = 1.006408249353399
$ws.Cells.Item(25, 8).Value = 1.023715350942155
$ws.Cells.Item(25, 10).Value = 0.08090115778746565
$ws.Cells.Item(25, 12).Value = 0.3724497180534456
$ws.Cells.Item(25, 13).Value = 0.3569387256756187
$ws.Cells.Item(25, 14).Value = 1.662131230317648
$ws.Cells.Item(25, 15).Value = 4.11395961754215
